# Refresh NATMI LR-pair edge-weight table (Bmp8a-Acvr2b) with newly
# recomputed TPM-based values. Only the "ECs" sending-cluster ligand
# expression (G/H) and "ECs" target-cluster receptor expression (M/N)
# changed at the source; every downstream specificity/edge-weight column
# (I, J, O, P, Q, R, S, T) is recomputed from those here to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.061724
$ws.Range("H2").Value = 0.185172
$ws.Range("I2").Value = 0.09652262708432048
$ws.Range("J2").Value = 0.09652262708432047
$ws.Range("M2").Value = 1.485259333333333
$ws.Range("N2").Value = 4.455778
$ws.Range("O2").Value = 0.3057455162066235
$ws.Range("P2").Value = 0.3057455162066235
$ws.Range("Q2").Value = 0.09167614709066668
$ws.Range("R2").Value = 0.8250853238160001
$ws.Range("S2").Value = 0.02951136044351498
$ws.Range("T2").Value = 0.02951136044351498
$ws.Range("G3").Value = 0.061724
$ws.Range("H3").Value = 0.185172
$ws.Range("I3").Value = 0.09652262708432048
$ws.Range("J3").Value = 0.09652262708432047
$ws.Range("O3").Value = 0.2805555239151429
$ws.Range("P3").Value = 0.2805555239151429
$ws.Range("Q3").Value = 0.08412306350933335
$ws.Range("R3").Value = 0.7571075715840001
$ws.Range("S3").Value = 0.0270799562113075
$ws.Range("T3").Value = 0.02707995621130749
$ws.Range("G4").Value = 0.061724
$ws.Range("H4").Value = 0.185172
$ws.Range("I4").Value = 0.09652262708432048
$ws.Range("J4").Value = 0.09652262708432047
$ws.Range("O4").Value = 0.4136989598782336
$ws.Range("P4").Value = 0.4136989598782336
$ws.Range("Q4").Value = 0.1240454060213333
$ws.Range("R4").Value = 1.116408654192
$ws.Range("S4").Value = 0.039931310429498
$ws.Range("T4").Value = 0.039931310429498
$ws.Range("I5").Value = 0.8735221647273214
$ws.Range("J5").Value = 0.8735221647273215
$ws.Range("M5").Value = 1.485259333333333
$ws.Range("N5").Value = 4.455778
$ws.Range("O5").Value = 0.3057455162066235
$ws.Range("P5").Value = 0.3057455162066235
$ws.Range("Q5").Value = 0.8296619029084444
$ws.Range("R5").Value = 7.466957126176
$ws.Range("S5").Value = 0.2670754851724821
$ws.Range("T5").Value = 0.2670754851724821
$ws.Range("I6").Value = 0.8735221647273214
$ws.Range("J6").Value = 0.8735221647273215
$ws.Range("O6").Value = 0.2805555239151429
$ws.Range("P6").Value = 0.2805555239151429
$ws.Range("S6").Value = 0.2450714685765634
$ws.Range("T6").Value = 0.2450714685765635
$ws.Range("I7").Value = 0.8735221647273214
$ws.Range("J7").Value = 0.8735221647273215
$ws.Range("O7").Value = 0.4136989598782336
$ws.Range("P7").Value = 0.4136989598782336
$ws.Range("S7").Value = 0.3613752109782759
$ws.Range("T7").Value = 0.361375210978276
$ws.Range("I8").Value = 0.02995520818835809
$ws.Range("J8").Value = 0.02995520818835809
$ws.Range("M8").Value = 1.485259333333333
$ws.Range("N8").Value = 4.455778
$ws.Range("O8").Value = 0.3057455162066235
$ws.Range("P8").Value = 0.3057455162066235
$ws.Range("Q8").Value = 0.02845113270288889
$ws.Range("R8").Value = 0.256060194326
$ws.Range("S8").Value = 0.00915867059062642
$ws.Range("T8").Value = 0.00915867059062642
$ws.Range("I9").Value = 0.02995520818835809
$ws.Range("J9").Value = 0.02995520818835809
$ws.Range("O9").Value = 0.2805555239151429
$ws.Range("P9").Value = 0.2805555239151429
$ws.Range("S9").Value = 0.008404099127271983
$ws.Range("T9").Value = 0.008404099127271983
$ws.Range("I10").Value = 0.02995520818835809
$ws.Range("J10").Value = 0.02995520818835809
$ws.Range("O10").Value = 0.4136989598782336
$ws.Range("P10").Value = 0.4136989598782336
$ws.Range("S10").Value = 0.01239243847045969
$ws.Range("T10").Value = 0.01239243847045969
